# Update the price list on "Hoja1" of the workbook:
#   - the list's issue date (A1) moves forward one month
#   - the six unit prices in column D get refreshed figures
#
# All six price cells keep their existing number format / style (0.000 via
# style index), so we only need to write the new numeric values - no
# formatting changes are required.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1: list date, stored as an Excel date-serial number (24-Apr-2024 -> 24-May-2024)
$ws.Range("A1").Value = 45436

# "SOPORTE LATERAL de Bce.1/2" / "...5/8"
$ws.Range("D22").Value = 400.797
$ws.Range("D23").Value = 508.443

# "SOPORTE TUBULAR de 1/2" / "...5/8"
$ws.Range("D34").Value = 396.131
$ws.Range("D35").Value = 548.628

# "SOPORTE CODO de Bce.1/2" / "...5/8"
$ws.Range("D45").Value = 492.295
$ws.Range("D46").Value = 545.4
